# ORG_HOSPITAL.xlsx edit
#
# The commit adds three new trailing columns (D, E, F) with header labels
# on row 1 of Sheet1, and moves the active selection to F5 (with the view
# scrolled so column B is the left-most visible column). No data is added
# to rows 2-12 for the new columns - only the headers are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (D1:F1). Setting these values is enough for
# Excel to grow the shared-string table, extend the sheet <dimension>,
# and widen every row's "spans" from "1:3" to "1:6" automatically.
$ws.Range("D1").Value = "ORG_HOS_IDENOLD"
$ws.Range("E1").Value = "ORG_HOS_IDENNEW"
$ws.Range("F1").Value = "ORG_HOS_STATUS"

# Move the selection to F5, matching the saved view state in the workbook.
[void]$ws.Range("F5").Select()
